# Auto-generated script to apply scheduled price-data refresh to Tonberry_Profits workbook
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for specific leve rows
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW worksheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 112
$ws.Range("H112").Value = 3763.1365
$ws.Range("J112").Value = 3763.1365
$ws.Range("L112").Value = 11289.4095
$ws.Range("N112").Value = -13505.4095
# Row 116
$ws.Range("H116").Value = 13698.728
$ws.Range("I116").Value = 34666.332
$ws.Range("K116").Value = 34666.332
$ws.Range("M116").Value = -31224.332
# Row 138
$ws.Range("H138").Value = 2205
$ws.Range("J138").Value = 2153.75
$ws.Range("L138").Value = 6461.25
$ws.Range("N138").Value = -16741.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 613.1429000000001
$ws.Range("I2").Value = 615.5
$ws.Range("K2").Value = 615.5
$ws.Range("M2").Value = -502.5
# Row 5
$ws.Range("H5").Value = 72
$ws.Range("I5").Value = 69.333336
$ws.Range("J5").Value = 80
$ws.Range("K5").Value = 69.333336
$ws.Range("L5").Value = 80
$ws.Range("M5").Value = 42.666664
$ws.Range("N5").Value = -304
# Row 74
$ws.Range("H74").Value = 4365
$ws.Range("I74").Value = 4352.7856
$ws.Range("K74").Value = 4352.7856
$ws.Range("M74").Value = -3478.7856
# Row 77
$ws.Range("H77").Value = 4365
$ws.Range("I77").Value = 4352.7856
$ws.Range("K77").Value = 21763.928
$ws.Range("M77").Value = -17395.928
# Row 116
$ws.Range("H116").Value = 613.1429000000001
$ws.Range("I116").Value = 615.5
$ws.Range("K116").Value = 615.5
$ws.Range("M116").Value = 1678.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 613.1429000000001
$ws.Range("I3").Value = 615.5
$ws.Range("K3").Value = 615.5
$ws.Range("M3").Value = -501.5
# Row 4
$ws.Range("H4").Value = 72
$ws.Range("I4").Value = 69.333336
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 69.333336
$ws.Range("L4").Value = 80
$ws.Range("M4").Value = 45.666664
$ws.Range("N4").Value = -310
# Row 20
$ws.Range("H20").Value = 1273.5172
$ws.Range("I20").Value = 1308
$ws.Range("J20").Value = 1196.8889
$ws.Range("K20").Value = 1308
$ws.Range("L20").Value = 1196.8889
$ws.Range("M20").Value = -1061
$ws.Range("N20").Value = -1690.8889
# Row 134
$ws.Range("H134").Value = 8283.556
$ws.Range("I134").Value = 9790.6
$ws.Range("K134").Value = 29371.8
$ws.Range("M134").Value = -26836.8

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2034
$ws.Range("I31").Value = 886.96295
$ws.Range("K31").Value = 886.96295
$ws.Range("M31").Value = -591.96295
# Row 34
$ws.Range("H34").Value = 2034
$ws.Range("I34").Value = 886.96295
$ws.Range("K34").Value = 886.96295
$ws.Range("M34").Value = -684.96295
# Row 62
$ws.Range("H62").Value = 3006.75
$ws.Range("I62").Value = 3244.75
$ws.Range("J62").Value = 2768.75
$ws.Range("K62").Value = 3244.75
$ws.Range("L62").Value = 2768.75
$ws.Range("M62").Value = -2620.75
$ws.Range("N62").Value = -4016.75
# Row 65
$ws.Range("H65").Value = 3006.75
$ws.Range("I65").Value = 3244.75
$ws.Range("J65").Value = 2768.75
$ws.Range("K65").Value = 16223.75
$ws.Range("L65").Value = 13843.75
$ws.Range("M65").Value = -13103.75
$ws.Range("N65").Value = -20083.75
# Row 68
$ws.Range("H68").Value = 44166.668
$ws.Range("J68").Value = 44166.668
$ws.Range("L68").Value = 44166.668
$ws.Range("N68").Value = -45664.668
# Row 71
$ws.Range("H71").Value = 44166.668
$ws.Range("J71").Value = 44166.668
$ws.Range("L71").Value = 132500.004
$ws.Range("N71").Value = -139988.004
# Row 74
$ws.Range("H74").Value = 30000
$ws.Range("J74").Value = 30000
$ws.Range("L74").Value = 30000
$ws.Range("N74").Value = -31748
# Row 77
$ws.Range("H77").Value = 30000
$ws.Range("J77").Value = 30000
$ws.Range("L77").Value = 90000
$ws.Range("N77").Value = -98736
# Row 93
$ws.Range("H93").Value = 50000
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 50000
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 50000
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -53744
# Row 107
$ws.Range("H107").Value = 387.17856
$ws.Range("I107").Value = 353.05
$ws.Range("J107").Value = 472.5
$ws.Range("K107").Value = 353.05
$ws.Range("L107").Value = 472.5
$ws.Range("M107").Value = 1566.95
$ws.Range("N107").Value = -4312.5
# Row 134
$ws.Range("H134").Value = 968.6
$ws.Range("I134").Value = 952.0714
$ws.Range("K134").Value = 2856.2142
$ws.Range("M134").Value = -321.2142000000003

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 498.75
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
# Row 122
$ws.Range("H122").Value = 1091.1538
$ws.Range("J122").Value = 1269.2858
$ws.Range("L122").Value = 11423.5722
$ws.Range("N122").Value = -16323.5722
# Row 131
$ws.Range("H131").Value = 19258308
$ws.Range("I131").Value = 71429144
$ws.Range("J131").Value = 37475
$ws.Range("K131").Value = 214287432
$ws.Range("L131").Value = 112425
$ws.Range("M131").Value = -214282392
$ws.Range("N131").Value = -122505
# Row 132
$ws.Range("H132").Value = 1890
$ws.Range("J132").Value = 2600
$ws.Range("L132").Value = 23400
$ws.Range("N132").Value = -28460
# Row 135
$ws.Range("H135").Value = 498.75
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
# Row 137
$ws.Range("H137").Value = 3394.842
$ws.Range("I137").Value = 1490.909
$ws.Range("J137").Value = 6012.75
$ws.Range("K137").Value = 4472.727000000001
$ws.Range("L137").Value = 18038.25
$ws.Range("M137").Value = 627.2729999999992
$ws.Range("N137").Value = -28238.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 22
$ws.Range("H22").Value = 55006
$ws.Range("I22").Value = 5000
$ws.Range("K22").Value = 5000
$ws.Range("M22").Value = -4471
# Row 122
$ws.Range("H122").Value = 1485.3478
$ws.Range("I122").Value = 1343.3125
$ws.Range("J122").Value = 1810
$ws.Range("K122").Value = 4029.9375
$ws.Range("L122").Value = 5430
$ws.Range("M122").Value = -1579.9375
$ws.Range("N122").Value = -10330

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1137.6
$ws.Range("I22").Value = 672
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 672
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -377
$ws.Range("N22").Value = -3590
# Row 27
$ws.Range("H27").Value = 1137.6
$ws.Range("I27").Value = 672
$ws.Range("J27").Value = 3000
$ws.Range("K27").Value = 672
$ws.Range("L27").Value = 3000
$ws.Range("M27").Value = -565
$ws.Range("N27").Value = -3214
# Row 122
$ws.Range("H122").Value = 5644.9546
$ws.Range("I122").Value = 2345
$ws.Range("J122").Value = 8394.916999999999
$ws.Range("K122").Value = 7035
$ws.Range("L122").Value = 25184.751
$ws.Range("M122").Value = -4585
$ws.Range("N122").Value = -30084.751

$wb.Save()
Write-Host "Applied scheduled price refresh to Tonberry_Profits workbook."